$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.3055
$ws.Range("E2").Value = -0.13485
$ws.Range("F2").Value = 0.0194
$ws.Range("G2").Value = 0.1501736054222039
$ws.Range("H2").Value = 0.1501736054222039
$ws.Range("I2").Value = 0.1233876453756994
$ws.Range("J2").Value = 0.1077179503346619
$ws.Range("K2").Value = 6183.8
$ws.Range("L2").Value = 0.08105256397317211
$ws.Range("M2").Value = 2647.3
$ws.Range("N2").Value = 0.01697886188456196
$ws.Range("O2").Value = 0.4281024612697694
$ws.Range("P2").Value = 2044.3
$ws.Range("Q2").Value = 0.01311142951331923
$ws.Range("R2").Value = 0.3305896050971894
$ws.Range("S2").Value = 603
$ws.Range("T2").Value = 0.2277792467797378
$ws.Range("U2").Value = 9336.200000000001
$ws.Range("V2").Value = 0.05987914113498558
$ws.Range("W2").Value = 0.06780475338614873
$ws.Range("X2").Value = 0.07539184375369752
$ws.Range("Y2").Value = -0.007587090367548788
$ws.Range("Z2").Value = 1.073007072867838
$ws.Range("AA2").Value = 0.09750796665805683
$ws.Range("AB2").Value = 0.06119126532424467
$ws.Range("AC2").Value = 0.03449238768610417
$ws.Range("AD2").Value = 20666.2
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 20666.2
$ws.Range("AG2").Value = 11330
$ws.Range("AH2").Value = 0.1170335184014824
$ws.Range("AI2").Value = 0.2234058193485339
$ws.Range("AJ2").Value = 0.06774395297027039
$ws.Range("AK2").Value = 0.1362286428837668
$ws.Range("AL2").Value = 626.4
$ws.Range("AM2").Value = 626.4
$ws.Range("AN2").Value = 2.100928156800553
$ws.Range("AO2").Value = 15.02825670498084
$ws.Range("AP2").Value = 1.151809041650147
$ws.Range("AQ2").Value = 15.02825670498084

# Row 3
$ws.Range("B3").Value = 'AIA Group Limited (SEHK:1299)'
$ws.Range("F3").Value = 0.0194
$ws.Range("G3").Value = 0.2161572559121818
$ws.Range("H3").Value = 0.2161572559121818
$ws.Range("I3").Value = 0.17326927986261
$ws.Range("J3").Value = 0.1493438183407055
$ws.Range("K3").Value = 5486
$ws.Range("L3").Value = 0.1273178769523544
$ws.Range("M3").Value = 1968
$ws.Range("N3").Value = 0.01331253475261549
$ws.Range("O3").Value = 0.3587313160772876
$ws.Range("P3").Value = 1965
$ws.Range("Q3").Value = 0.01329224125451699
$ws.Range("R3").Value = 0.3581844695588771
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 0.001524390243902439
$ws.Range("U3").Value = 5027
$ws.Range("V3").Value = 0.03400513831371854
$ws.Range("W3").Value = 0.1099883716267693
$ws.Range("X3").Value = 0.06319627948550255
$ws.Range("Y3").Value = 0.04679209214126677
$ws.Range("Z3").Value = 0.7740770681756939
$ws.Range("AA3").Value = 0.1156036250513368
$ws.Range("AB3").Value = 0.06073929718316007
$ws.Range("AC3").Value = 0.05486432786817669
$ws.Range("AD3").Value = 9394
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 9394
$ws.Range("AG3").Value = 4367
$ws.Range("AH3").Value = 0.05974891969831693
$ws.Range("AI3").Value = 0.1409473510480277
$ws.Range("AJ3").Value = 0.02869296230689577
$ws.Range("AK3").Value = 0.0708675473045341
$ws.Range("AL3").Value = 290
$ws.Range("AM3").Value = 290
$ws.Range("AN3").Value = 1.226530878704792
$ws.Range("AO3").Value = 25.7448275862069
$ws.Range("AP3").Value = 0.5701788745267006
$ws.Range("AQ3").Value = 25.7448275862069
$ws.Range("D3").Value = $null

# Row 4
$ws.Range("B4").Value = 'Yunfeng Financial Group Limited (SEHK:376)'
$ws.Range("D4").Value = 0.479
$ws.Range("E4").Value = -0.186
$ws.Range("G4").Value = 0.09403058707449433
$ws.Range("H4").Value = 0.09403058707449433
$ws.Range("I4").Value = 0.1003453379378392
$ws.Range("J4").Value = 0.1001114411960276
$ws.Range("K4").Value = 34.5
$ws.Range("L4").Value = 0.03404045387271831
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 307.4
$ws.Range("V4").Value = 0.191359561752988
$ws.Range("W4").Value = 0.02872127872127872
$ws.Range("X4").Value = 0.07539184375369752
$ws.Range("Y4").Value = -0.0466705650324188
$ws.Range("Z4").Value = 0.9557714070162202
$ws.Range("AA4").Value = 0.09568365301034884
$ws.Range("AB4").Value = 0.06119126532424467
$ws.Range("AC4").Value = 0.03449238768610417
$ws.Range("AD4").Value = 560.5
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 560.5
$ws.Range("AG4").Value = 253.1
$ws.Range("AH4").Value = 0.2586644515206055
$ws.Range("AI4").Value = 0.2160838891244844
$ws.Range("AJ4").Value = 0.1361118580263512
$ws.Range("AK4").Value = 0.1106931992127706
$ws.Range("AL4").Value = 18.7
$ws.Range("AM4").Value = 18.7
$ws.Range("AN4").Value = 5
$ws.Range("AO4").Value = 5.438502673796791
$ws.Range("AP4").Value = 2.257805530776093
$ws.Range("AQ4").Value = 5.438502673796791
$ws.Range("T4").Value = $null

# Row 5
$ws.Range("D5").Value = 0.132
$ws.Range("E5").Value = -0.0837
$ws.Range("G5").Value = 0.06361987126916673
$ws.Range("H5").Value = 0.06361987126916673
$ws.Range("I5").Value = 0.05734486443500087
$ws.Range("J5").Value = 0.0435492087599369
$ws.Range("K5").Value = 663.3
$ws.Range("L5").Value = 0.02060501006486245
$ws.Range("M5").Value = 679.3
$ws.Range("N5").Value = 0.1048237763101043
$ws.Range("O5").Value = 1.024121815166591
$ws.Range("P5").Value = 79.3
$ws.Range("Q5").Value = 0.01223689895685452
$ws.Range("R5").Value = 0.1195537464194181
$ws.Range("S5").Value = 600
$ws.Range("T5").Value = 0.8832621816575887
$ws.Range("U5").Value = 4001.8
$ws.Range("V5").Value = 0.6175236096537251
$ws.Range("W5").Value = 0.06780475338614873
$ws.Range("X5").Value = 0.1311202396068598
$ws.Range("Y5").Value = -0.06331548622071108
$ws.Range("Z5").Value = 2.23902958135394
$ws.Range("AA5").Value = 0.09750796665805683
$ws.Range("AB5").Value = 0.0648457138682694
$ws.Range("AC5").Value = 0.03266225278978743
$ws.Range("AD5").Value = 10711.7
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 10711.7
$ws.Range("AG5").Value = 6709.900000000001
$ws.Range("AH5").Value = 0.62305942845842
$ws.Range("AI5").Value = 0.4604746736135292
$ws.Range("AJ5").Value = 0.5086995746874597
$ws.Range("AK5").Value = 0.3483762103787545
$ws.Range("AL5").Value = 317.7
$ws.Range("AM5").Value = 317.7
$ws.Range("AN5").Value = 5.185757164988382
$ws.Range("AO5").Value = 5.810513062637709
$ws.Range("AP5").Value = 3.24840240123935
$ws.Range("AQ5").Value = 5.810513062637709
